$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 187, shifting existing rows 187:219 down to 188:220
$ws.Rows.Item(187).EntireRow.Insert()

# Populate the newly inserted row 187 with the new record
$ws.Range("A187").Value = 8
$ws.Range("B187").Value = "Terminal La Palmera de La Serena"
$ws.Range("C187").Value = "Coquimbo"
$ws.Range("D187").Value = 45015
$ws.Range("E187").Value = 4
$ws.Range("F187").Value = 100112001
$ws.Range("G187").Value = "Berenjena"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 500
$ws.Range("K187").Value = 10500
$ws.Range("L187").Value = 11000
$ws.Range("M187").Value = 10750
$ws.Range("N187").Value = "$/caja 50 unidades"
$ws.Range("O187").Value = "Región de Arica y Parinacota"
$ws.Range("P187").Value = 215
$ws.Range("Q187").Value = 50
$ws.Range("R187").Value = "Hortaliza"
